$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 314
$ws.Range("C3").Value = "naman roy"
$ws.Range("D3").Value = "roy1998@gmail.com"
$ws.Range("G3").Value = "L1_selected"
